$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "SCD0016"

# 2. Update column B content (ticket id) for all data rows 2-7
$ws.Range("B2:B7").Value = "SCD0016-041"

# 3. Apply formatting to the used data range A1:P7
$fmtRange = $ws.Range("A1:P7")
$fmtRange.Font.Size = 10
$fmtRange.HorizontalAlignment = -4131  # xlLeft
$fmtRange.VerticalAlignment = -4108    # xlCenter

# 4. Restore selection/scroll position
$ws.Range("A5").Select()
$ws.Range("B8").Select()
